$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.578.37"
$ws.Range("E2").Value = "  +0.75%  "

$ws.Range("D3").Value = "1.876.45"
$ws.Range("E3").Value = "  -0.07%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.0000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.44%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "248.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.23%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9999"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.39%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4760"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2910"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.99%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06503"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.23%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.95"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.76%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07747"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.38%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7379"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.24%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "96.33"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.17%  "

$ws.Range("D14").Value = "1.873.19"
$ws.Range("E14").Value = "  -0.31%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.174"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.66%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "274.31"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.15%  "

$ws.Range("D17").Value = "30.644.43"
$ws.Range("E17").Value = "  +0.89%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.22"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.34%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.0000"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.26%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007513"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.25%  "

$ws.Range("D21").Value = "2.118.99"
$ws.Range("E21").Value = "  -0.59%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.0000"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.53%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.242"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.06%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.177"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.06%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.189"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.78%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.91"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.02%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.81"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.85%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.908"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.04%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.09850"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.75%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.339"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.57%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.499"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.23%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.265"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.07%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.088"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.01%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04810"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.25%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.122"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.04%  "

$ws.Range("E36").Value = "  -0.41%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.717"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.41%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01860"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.14%  "

$ws.Range("E39").Value = "  +0.54%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.283"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.89%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "73.42"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.49%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.984"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.02%  "

$ws.Range("E43").Value = "  +0.73%  "

$ws.Range("E44").Value = "  -0.35%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8349"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.16%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "101.73"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.21%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.388"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.87%  "

$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.982"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.75%  "

$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "35.30"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.48%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "914.18"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.89%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05673"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.43%  "
